# Apply the changes described by the commit:
#  - Re-brand URL from ibm.com to linuxforhealth.org
#  - Bump version 7.0.0 -> 8.0.0
#  - Update the publication Date
#  - Rename Publisher from "Alvearie Team" to "LinuxForHealth Team"
#  - Clear the Constraint(s) cell on the root "Extension" row of the
#    Elements sheet (the constraint now only applies to Extension.extension)

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-union-id"
$wsMetadata.Range("B3").Value = "8.0.0"
$wsMetadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMetadata.Range("B9").Value = "LinuxForHealth Team"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("AI2").Value = ""
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-union-id"
